# Refresh the "data pobrania" (download timestamp) column (V) for every
# data row with the new scrape time recorded for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-01-05 22:42:05"

$lastRow = $ws.UsedRange.Rows.Count
$firstDataRow = 2

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 22).Value = $newTimestamp
}
